# Apply "hybrid bold + color highlighting" to quantitative metrics in
# specific bullet/paragraph lines, matching the target OOXML diff.
#
# Strategy: for each target paragraph (identified by a unique anchor
# substring, optionally combined with an "exclude" substring to
# disambiguate from a near-duplicate paragraph), locate the paragraph,
# then within that paragraph's range use Find scoped left-to-right to
# locate each metric substring in turn and apply Bold + the highlight
# color (#2C3E50) to just that substring. Word's own run-splitting takes
# care of turning the single run into the prefix/bold/suffix runs seen
# in the diff.
#
# NOTE: this interpreter's function calls only reliably bind POSITIONAL
# parameters (named `-param value` binding does not populate $param), so
# all helper-function calls below use positional args only.

$d = $word.ActiveDocument
$highlightColor = 5258796   # OLE BGR decimal for hex 2C3E50 (R=44,G=62,B=80)

function Highlight-MetricsInParagraph {
    param(
        [string]$anchorText,    # substring identifying the paragraph
        [string]$excludeText,   # if non-empty, paragraphs containing this are skipped (disambiguation)
        [string[]]$metrics      # ordered list of substrings to bold, left-to-right
    )

    $target = $null
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        $pText = $p.Range.Text
        if ($pText -like "*$anchorText*") {
            if ($excludeText -ne "" -and $pText -like "*$excludeText*") {
                continue
            }
            $target = $p
            break
        }
    }

    if ($null -eq $target) {
        Write-Host "WARNING: paragraph not found for anchor: $anchorText"
        return
    }

    $pStart = $target.Range.Start
    $pEnd = $target.Range.End
    $searchFrom = $pStart

    foreach ($metric in $metrics) {
        $r = $d.Range($searchFrom, $pEnd)
        $found = $r.Find.Execute($metric, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if ($found) {
            $r.Font.Bold = 1
            $r.Font.Color = $highlightColor
            $searchFrom = $r.End
        } else {
            Write-Host "WARNING: metric not found in paragraph: $metric"
        }
    }
}

# 1. "Discovered systematic race coding errors ... from 23% to 64%"
Highlight-MetricsInParagraph "Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms" "" @("23%", "64%")

# 2. "Achieved 87% prediction accuracy ... from ±4.2% to ±2.1%" (long version, with polling error margins)
$pm = [char]0x00B1
Highlight-MetricsInParagraph "reducing polling error margins" "" @("87%", "71%", ($pm + "4.2%"), ($pm + "2.1%"))

# 3. "Wrote RFP and analyzed bids from 1,200 vendors ..."
Highlight-MetricsInParagraph "Wrote RFP and analyzed bids from" "" @("1,200")

# 4. "Created comprehensive meta-analysis framework ... $400M ... $1B+"
Highlight-MetricsInParagraph "Created comprehensive meta-analysis framework" "" @('$400M', '$1B')

# 5. "Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M"
Highlight-MetricsInParagraph "Algorithm reduced mapping costs by" "" @("73.5%", '$4.7M')

# 6. "Achieved 87% prediction accuracy ... industry standard of 71%" (short version, in KEY ACHIEVEMENTS)
# Exclude the long "reducing polling error margins" paragraph to disambiguate.
Highlight-MetricsInParagraph "Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%" "reducing polling error margins" @("87%", "71%")
